$d = $word.ActiveDocument

# 1. Title: "from summer 2019 " -> "from 2019 "
$d.Content.Find.Execute(
    "Optical, Chemical, and Biological Oceanographic Conditions in the Labrador Sea from summer 2019 ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Optical, Chemical, and Biological Oceanographic Conditions in the Labrador Sea from 2019 ",
    2) | Out-Null

# 2. "Ocean and Ecosystem Sciences Division, Science Branch" -> "Fisheries and Oceans Canada"
$d.Content.Find.Execute(
    "Ocean and Ecosystem Sciences Division, Science Branch",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Fisheries and Oceans Canada",
    2) | Out-Null

# 3. "Department of Fisheries and Oceans" -> "Science Branch, Maritimes Region"
$d.Content.Find.Execute(
    "Department of Fisheries and Oceans",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Science Branch, Maritimes Region",
    2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "$i : [$($p.Range.Text)]"
}
